$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the project title in A11 from "Blur Buddy" to "Censor Buddy"
$ws.Range("A11").Value = "Censor Buddy"

# Move the active cell selection to A12 (as left after the edit)
$ws.Range("A12").Select()
